$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "79174445"
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = "Cash"
$ws.Range("H2").Value = "2025-08-21T07:35:23"
